# Auto-generated edit script: update price/profit figures per the commit diff
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 874.1515000000001
$ws.Range("I15").Value = 874.1515000000001
$ws.Range("K15").Value = 2622.4545
$ws.Range("M15").Value = -2453.4545
$ws.Range("H32").Value = 2168.2856
$ws.Range("J32").Value = 1933.5
$ws.Range("L32").Value = 1933.5
$ws.Range("N32").Value = -2585.5
$ws.Range("H53").Value = 35103.965
$ws.Range("I53").Value = 787.64703
$ws.Range("J53").Value = 83718.75
$ws.Range("K53").Value = 787.64703
$ws.Range("L53").Value = 83718.75
$ws.Range("M53").Value = -150.64703
$ws.Range("N53").Value = -84992.75
$ws.Range("H98").Value = 3078.1936
$ws.Range("I98").Value = 2475.9048
$ws.Range("J98").Value = 4343
$ws.Range("K98").Value = 2475.9048
$ws.Range("L98").Value = 4343
$ws.Range("M98").Value = -977.9047999999998
$ws.Range("N98").Value = -7339
$ws.Range("H122").Value = 3078.1936
$ws.Range("I122").Value = 2475.9048
$ws.Range("J122").Value = 4343
$ws.Range("K122").Value = 7427.714399999999
$ws.Range("L122").Value = 13029
$ws.Range("M122").Value = -4977.714399999999
$ws.Range("N122").Value = -17929
$ws.Range("H123").Value = 62687.5
$ws.Range("J123").Value = 73583.336
$ws.Range("L123").Value = 73583.336
$ws.Range("N123").Value = -83383.336
$ws.Range("H132").Value = 4315.2607
$ws.Range("I132").Value = 4422.45
$ws.Range("K132").Value = 13267.35
$ws.Range("M132").Value = -10737.35
$ws.Range("H137").Value = 7030.3076
$ws.Range("J137").Value = 7485.143
$ws.Range("L137").Value = 22455.429
$ws.Range("N137").Value = -27555.429
$ws.Range("H138").Value = 6888.702
$ws.Range("I138").Value = 6367.091
$ws.Range("J138").Value = 7048.0835
$ws.Range("K138").Value = 19101.273
$ws.Range("L138").Value = 21144.2505
$ws.Range("M138").Value = -13961.273
$ws.Range("N138").Value = -31424.2505

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4431
$ws.Range("I5").Value = 303.4
$ws.Range("J5").Value = 14750
$ws.Range("K5").Value = 303.4
$ws.Range("L5").Value = 14750
$ws.Range("M5").Value = -191.4
$ws.Range("N5").Value = -14974
$ws.Range("H11").Value = 6668667.5
$ws.Range("I11").Value = 6668667.5
$ws.Range("K11").Value = 6668667.5
$ws.Range("M11").Value = -6668523.5
$ws.Range("H21").Value = 1500
$ws.Range("I21").Value = 1500
$ws.Range("K21").Value = 1500
$ws.Range("M21").Value = -1126
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10804
$ws.Range("H45").Value = 310387.84
$ws.Range("I45").Value = 372144.47
$ws.Range("J45").Value = 1604.6666
$ws.Range("K45").Value = 372144.47
$ws.Range("L45").Value = 1604.6666
$ws.Range("M45").Value = -371767.47
$ws.Range("N45").Value = -2358.6666
$ws.Range("H74").Value = 2533.186
$ws.Range("I74").Value = 1965.0358
$ws.Range("J74").Value = 3593.7334
$ws.Range("K74").Value = 1965.0358
$ws.Range("L74").Value = 3593.7334
$ws.Range("M74").Value = -1091.0358
$ws.Range("N74").Value = -5341.7334
$ws.Range("H77").Value = 2533.186
$ws.Range("I77").Value = 1965.0358
$ws.Range("J77").Value = 3593.7334
$ws.Range("K77").Value = 9825.179
$ws.Range("L77").Value = 17968.667
$ws.Range("M77").Value = -5457.179
$ws.Range("N77").Value = -26704.667
$ws.Range("H131").Value = 99998
$ws.Range("J131").Value = 99998
$ws.Range("L131").Value = 99998
$ws.Range("N131").Value = -110078
$ws.Range("H132").Value = 14447188
$ws.Range("I132").Value = 15456750
$ws.Range("K132").Value = 46370250
$ws.Range("M132").Value = -46367720

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4431
$ws.Range("I4").Value = 303.4
$ws.Range("J4").Value = 14750
$ws.Range("K4").Value = 303.4
$ws.Range("L4").Value = 14750
$ws.Range("M4").Value = -188.4
$ws.Range("N4").Value = -14980
$ws.Range("H107").Value = 20836088
$ws.Range("I107").Value = 2763.7222
$ws.Range("J107").Value = 83336056
$ws.Range("K107").Value = 2763.7222
$ws.Range("L107").Value = 83336056
$ws.Range("M107").Value = -843.7222000000002
$ws.Range("N107").Value = -83339896
$ws.Range("H134").Value = 4284.913
$ws.Range("I134").Value = 3932.9412
$ws.Range("K134").Value = 11798.8236
$ws.Range("M134").Value = -9263.8236

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 91.888885
$ws.Range("I7").Value = 76.933334
$ws.Range("J7").Value = 166.66667
$ws.Range("K7").Value = 76.933334
$ws.Range("L7").Value = 166.66667
$ws.Range("M7").Value = 36.066666
$ws.Range("N7").Value = -392.66667
$ws.Range("H13").Value = 400
$ws.Range("I13").Value = 400
$ws.Range("K13").Value = 400
$ws.Range("M13").Value = -261
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H105").Value = 852.3333
$ws.Range("I105").Value = 848.0909
$ws.Range("K105").Value = 848.0909
$ws.Range("M105").Value = 898.9091
$ws.Range("H132").Value = 628263.9
$ws.Range("I132").Value = 3053.875
$ws.Range("J132").Value = 1253473.9
$ws.Range("K132").Value = 9161.625
$ws.Range("L132").Value = 3760421.7
$ws.Range("M132").Value = -6631.625
$ws.Range("N132").Value = -3765481.7
$ws.Range("H141").Value = 1055999.6
$ws.Range("J141").Value = 1055999.6
$ws.Range("L141").Value = 1055999.6
$ws.Range("N141").Value = -1066359.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4224.5
$ws.Range("J39").Value = 4224.5
$ws.Range("L39").Value = 12673.5
$ws.Range("N39").Value = -13261.5
$ws.Range("H55").Value = 2049.5
$ws.Range("J55").Value = 2466
$ws.Range("L55").Value = 7398
$ws.Range("N55").Value = -7752
$ws.Range("H56").Value = 5166.6665
$ws.Range("I56").Value = 5166.6665
$ws.Range("K56").Value = 5166.6665
$ws.Range("M56").Value = -4636.6665
$ws.Range("H122").Value = 1082923.1
$ws.Range("I122").Value = 726.5
$ws.Range("J122").Value = 6493906
$ws.Range("K122").Value = 6538.5
$ws.Range("L122").Value = 58445154
$ws.Range("M122").Value = -4088.5
$ws.Range("N122").Value = -58450054

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4999
$ws.Range("I4").Value = 4999
$ws.Range("M4").Value = -4887
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83338170
$ws.Range("J7").Value = 6249.5
$ws.Range("L7").Value = 6249.5
$ws.Range("N7").Value = -6473.5
$ws.Range("H22").Value = 1590
$ws.Range("J22").Value = 2333.3333
$ws.Range("L22").Value = 2333.3333
$ws.Range("N22").Value = -2923.3333
$ws.Range("H27").Value = 1590
$ws.Range("J27").Value = 2333.3333
$ws.Range("L27").Value = 2333.3333
$ws.Range("N27").Value = -2547.3333
$ws.Range("H36").Value = 61000
$ws.Range("J36").Value = 61000
$ws.Range("L36").Value = 61000
$ws.Range("N36").Value = -62124
$ws.Range("H46").Value = 348042.12
$ws.Range("I46").Value = 1323
$ws.Range("J46").Value = 403517.2
$ws.Range("K46").Value = 1323
$ws.Range("L46").Value = 403517.2
$ws.Range("M46").Value = -1135
$ws.Range("N46").Value = -403893.2
$ws.Range("H126").Value = 83338170
$ws.Range("J126").Value = 6249.5
$ws.Range("L126").Value = 18748.5
$ws.Range("N126").Value = -23688.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H43").Value = 47100.36
$ws.Range("I43").Value = 5109.1304
$ws.Range("K43").Value = 5109.1304
$ws.Range("M43").Value = -4960.1304
$ws.Range("H47").Value = 36227.6
$ws.Range("J47").Value = 36227.6
$ws.Range("L47").Value = 36227.6
$ws.Range("N47").Value = -37371.6
$ws.Range("H48").Value = 10010000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21138
$ws.Range("H81").Value = 7412698.5
$ws.Range("I81").Value = 4563.1665
$ws.Range("J81").Value = 22228970
$ws.Range("K81").Value = 9126.333000000001
$ws.Range("L81").Value = 44457940
$ws.Range("M81").Value = -8065.333000000001
$ws.Range("N81").Value = -44460062
$ws.Range("H84").Value = 7412698.5
$ws.Range("I84").Value = 4563.1665
$ws.Range("J84").Value = 22228970
$ws.Range("K84").Value = 45631.665
$ws.Range("L84").Value = 222289700
$ws.Range("M84").Value = -40327.665
$ws.Range("N84").Value = -222300308
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 20000
$ws.Range("K88").Value = 20000
$ws.Range("M88").Value = -19594
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 20000
$ws.Range("K91").Value = 20000
$ws.Range("M91").Value = -18596
$ws.Range("H126").Value = 5145.636
$ws.Range("I126").Value = 2857.7144
$ws.Range("K126").Value = 8573.143199999999
$ws.Range("M126").Value = -6103.143199999999
$ws.Range("H132").Value = 274789.66
$ws.Range("I132").Value = 290106.2
$ws.Range("K132").Value = 870318.6000000001
$ws.Range("M132").Value = -867788.6000000001
$ws.Range("H136").Value = 8801.839
$ws.Range("I136").Value = 9124.406999999999
$ws.Range("K136").Value = 27373.221
$ws.Range("M136").Value = -24823.221
